# Update item list: add three new products (Mini poster, Photobook, Star
# Keychain), which Excel's alphabetical sort on column B (name) interleaves
# into the existing rows. Rows 21-36 hold the new sorted tail of the table;
# rows 2-20 are unaffected by the re-sort. Re-sorting + appending in the UI
# produced this exact layout, so we just rewrite the affected cells directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(35, 'Mini poster', 15000, 'Prints', 'Artboard 36.png', 0),
    @(7, 'No way', 9000, 'Sticker', 'Artboard 8.png', 0),
    @(33, 'Photobook', 120000, 'Other', 'Artboard 36.png', 0),
    @(30, 'Playdate a4', 50000, 'Prints', 'Artboard 33.png', 0),
    @(31, 'Pouch ', 30000, 'Other', 'Artboard 34.png', 0),
    @(16, 'Prot prot', 8000, 'Sticker', 'Artboard 18.png', 0),
    @(14, 'Safe place custard', 9000, 'Sticker', 'Artboard 15.png', 0),
    @(10, 'Sleepy joy', 11000, 'Sticker', 'Artboard 11.png', 0),
    @(34, 'Star Keychain', 75000, 'Keychain', 'Artboard 35.png', 0),
    @(19, 'Starry Custard', 9000, 'Sticker', 'Artboard 21.png', 0),
    @(8, 'Starry nori', 8000, 'Sticker', 'Artboard 9.png', 0),
    @(32, 'Studying churros', 9000, 'Sticker', 'Artboard 16.png', 0),
    @(29, 'Taiyaki girl', 65000, 'Keychain', 'Artboard 32.png', 0),
    @(15, 'Unicycle joy', 9000, 'Sticker', 'Artboard 17.png', 0),
    @(18, 'Worried churros', 9000, 'Sticker', 'Artboard 20.png', 0),
    @(9, 'Zooming churros', 11000, 'Sticker', 'Artboard 10.png', 0)
)

$r = 21
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Selection ends up on the newly-added last row after the insert/sort.
$ws.Range("B36").Select()
